$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.731.28'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.26%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.872.47'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.00%  '

$ws.Range('E4').Value = '  +0.41%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.57'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.03%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.31%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4624'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.30%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3855'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.38%  '

$ws.Range('E9').Value = '  +0.10%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9764'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.60%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.82'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.51%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.847.03'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.45%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.015'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.46%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.705'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.28%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06957'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.63%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.51'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.42%  '

$ws.Range('E17').Value = '  +0.37%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001004'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.14%  '

$ws.Range('E19').Value = '  +0.96%  '

$ws.Range('E20').Value = '  +0.35%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '28.716.48'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.21%  '

$ws.Range('E22').Value = '  -1.03%  '

$ws.Range('E23').Value = '  +0.67%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.098'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.24%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.064.72'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.90%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.08'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.72%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.26'
$ws.Range('D27').Style = 'Normal'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.885'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.22%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.986'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.00%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '119.18'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.94%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09334'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.86%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9181'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.00%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.292'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.38%  '

$ws.Range('E34').Value = '  +0.85%  '

$ws.Range('E35').Value = '  +0.86%  '

$ws.Range('E36').Value = '  -0.91%  '

$ws.Range('E37').Value = '  +0.95%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02078'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.35%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.658'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.52%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5627'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.52%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1784'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.09%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.776'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.28%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.07224'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.18%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.71'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.74%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5291'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.56%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.147'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.85%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.120'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.18%  '

$ws.Range('E48').Value = '  +0.09%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '112.94'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.20%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.410'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.73%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.003'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.35%  '
